$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("RegNum", "Make", "Colour", "Result"),
    @("EF63 YPZ", "FORD", "RED", "Pass"),
    @("KY66 LNU", "SMART", "BLACK", "Pass"),
    @("PE65 YNY", "RENAULT", "WHITE", "Pass"),
    @("YC11 OMK", "MINI", "SILVER", "Pass"),
    @("UK65 FCV", "TOYOTA", "BLUE", "Pass"),
    @("KR66 NUO", "SUZUKI", "BLUE", "Pass")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Page setup: paper size + orientation (print settings)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("D14").Select()
